$d = $word.ActiveDocument

$pairs = @(
    @{old="80×28="; new="14×68="},
    @{old="60×61="; new="73×46="},
    @{old="77×28="; new="68×59="},
    @{old="96×26="; new="25×58="},
    @{old="64×22="; new="88×91="},
    @{old="72×99="; new="27×34="},
    @{old="76×63="; new="63×80="},
    @{old="78×46="; new="28×32="},
    @{old="72×88="; new="68×17="},
    @{old="41×76="; new="41×61="},
    @{old="38×53="; new="69×45="},
    @{old="34×13="; new="61×77="},
    @{old="46×56="; new="47×37="},
    @{old="55×53="; new="32×11="},
    @{old="20×62="; new="11×48="},
    @{old="65×34="; new="52×47="},
    @{old="94×97="; new="92×66="},
    @{old="73×29="; new="87×57="},
    @{old="56×49="; new="65×80="},
    @{old="69×34="; new="47×99="},
    @{old="26×25="; new="36×27="},
    @{old="29×61="; new="32×62="},
    @{old="86×26="; new="47×90="},
    @{old="54×98="; new="37×77="},
    @{old="44×94="; new="52×39="}
)

foreach ($p in $pairs) {
    $range = $d.Content
    $range.Find.Execute($p.old, $true, $false, $false, $false, $false, $true, 1, $false, $p.new, 2)
}
